# Update input data with electric boiler info.
# - Row 7 on the "Units" sheet (Steam_Plant) is re-pointed from a
#   Water+Waste_Heat fed unit to a Power_Kasso (electricity) + Water fed
#   electric boiler, with refreshed capacity / efficiency figures.
# - X5 (Relation_In_In for the Destilation_Tower unit) gets a corrected value.
# - Selection/active-sheet view state is updated: "Units" becomes the
#   selected/active tab (instead of "Connections"), with a new selection.

$wb = $excel.ActiveWorkbook

$units = $wb.Worksheets.Item("Units")

# --- Units!X5 : Relation_In_In value correction -----------------------------
$units.Range("X5").Value = 0.0079901515151515144

# --- Units row 7 (Steam_Plant -> electric boiler) ---------------------------
# Input1: Water -> Power_Kasso
$units.Range("B7").Value = "Power_Kasso"
# Input2: Waste_Heat -> Water
$units.Range("C7").Value = "Water"
# Output1 (Steam) is unchanged.

# Cap_Output1_existing (J7) is no longer used; Cap_Input1_existing (F7) is set instead.
$units.Range("J7").Value = $null
$units.Range("F7").Value = 100

# Relation_In_In
$units.Range("X7").Value = 1
# Relation_In_Out (new)
$units.Range("Y7").Value = 0.99
# fom_cost (new)
$units.Range("AC7").Value = 0.11929223744292237
# unit_idle_heat_rate (new)
$units.Range("AJ7").Value = 0.1

# --- View state: make "Units" the active/selected sheet ---------------------
$units.Activate()
$units.Range("X8").Select()

Write-Output "electric boiler info applied"
